$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Purpose" paragraph (paragraph 3) - split the sentence so it talks
# about administrators backing up / restoring the CMS database instead of
# users signing in via SSO.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$xml1 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5DBE1E2C" w14:textId="7F7DFDE2" w:rsidR="00EF32E2" w:rsidRDefault="00EF32E2" w:rsidP="00EF32E2">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="00EF32E2">
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve">To </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve">allow for </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve">administrators to </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    </w:rPr>
    <w:t>backup</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve"> the CMS database and restore.</w:t>
  </w:r>
</w:p>
"@
$p3.Range.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------------
# Edit 2: "Background" paragraph - trim the long sentence down to just the
# first clause, then add a new second run about 1Gov's redundancy policy /
# administrator training.
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$xml2 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2ED52FA6" w14:textId="0173EFCE" w:rsidR="00EF32E2" w:rsidRDefault="00A80880" w:rsidP="00EF32E2">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve">The MESD LMS will be hosted and housed by the 1Gov 1Citizen platform. </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    </w:rPr>
    <w:t>Because of this it will be up to 1Gov&#8217;s redundancy policy to know how often and how much backups of the database will be available. Administrators will be trained to get necessary support during the TOT workshop.</w:t>
  </w:r>
</w:p>
"@
$p6.Range.InsertXML($xml2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 3: remove the "SSO is accomplished..." paragraph's text plus every
# paragraph that follows it (Redirect Module, plugin descriptions, setup
# process, screenshots, SSO related plugins, ...) leaving a single empty
# paragraph where that content used to start.
# ---------------------------------------------------------------------------
$startPara = $d.Paragraphs.Item(7)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng = $d.Range($startPara.Range.Start, $lastPara.Range.End)
$xml3 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1AD56758" w14:textId="44F3F4D9" w:rsidR="00245634" w:rsidRDefault="00245634" w:rsidP="00EF32E2">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@
$rng.InsertXML($xml3) | Out-Null
